$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with the results of the new test run (28 jun 2023 execution)
$ws.Range("D2").Value = "/src/Excel/entregable2/ABONOCTS2023020101.txt"
$ws.Range("E2").Value = "PASSED"
$ws.Range("F2").Value = "UPLD2317909159 14"
$ws.Range("G2").Value = "28 jun. 2023, 14:35:52"

# Column D auto-narrows (best-fit) since the new file path is shorter than before
$ws.Columns.Item(4).ColumnWidth = 48.8

# Cursor/selection position as last left in the workbook view
$ws.Range("F11").Select() | Out-Null
